$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue "D2" "261.33"
Set-TextValue "E2" "1.70%"
Set-TextValue "D3" "27.36"
Set-TextValue "E3" "1.27%"
Set-TextValue "D4" "4.710"
Set-TextValue "E4" "3.53%"
Set-TextValue "D5" "0.06082"
Set-TextValue "E5" "3.16%"
Set-TextValue "D6" "6.679"
Set-TextValue "E6" "1.06%"
Set-TextValue "D7" "0.8456"
Set-TextValue "E7" "-0.51%"
Set-TextValue "D8" "0.9234"
Set-TextValue "E8" "-0.80%"
Set-TextValue "D9" "0.1403"
Set-TextValue "E9" "2.04%"
Set-TextValue "D10" "0.05025"
Set-TextValue "E10" "19.19%"
Set-TextValue "D11" "0.07141"
Set-TextValue "E11" "1.84%"
Set-TextValue "D12" "0.03126"
Set-TextValue "E12" "2.65%"
Set-TextValue "D13" "0.09070"
Set-TextValue "E13" "-0.38%"
Set-TextValue "D14" "0.001546"
Set-TextValue "E14" "1.26%"
Set-TextValue "D15" "0.0006079"
Set-TextValue "E15" "0.70%"
Set-TextValue "D16" "0.006074"
Set-TextValue "E16" "-0.14%"
Set-TextValue "E17" "-0.49%"
Set-TextValue "E18" "-0.80%"
Set-TextValue "D19" "2.168"
Set-TextValue "E19" "-2.00%"
Set-TextValue "D21" "0.1306"
Set-TextValue "E21" "0.89%"
Set-TextValue "D22" "4.091"
Set-TextValue "E22" "4.74%"
Set-TextValue "D23" "0.04257"
Set-TextValue "E23" "-0.09%"
Set-TextValue "D24" "0.001221"
Set-TextValue "E24" "-0.03%"
Set-TextValue "D26" "0.0001200"
Set-TextValue "E26" "0.06%"
Set-TextValue "E27" "3.47%"
Set-TextValue "E40" "1.97%"
Set-TextValue "E41" "1.43%"
Set-TextValue "D42" "0.004093"
Set-TextValue "E42" "-35.01%"
Set-TextValue "D43" "0.01636"
Set-TextValue "E43" "21.65%"
Set-TextValue "E44" "1.00%"
Set-TextValue "D45" "0.00005260"
Set-TextValue "E45" "-1.76%"
Set-TextValue "E46" "0.10%"
Set-TextValue "E47" "6.96%"
Set-TextValue "D48" "0.1358"
Set-TextValue "E48" "-46.12%"
Set-TextValue "E49" "0.10%"
Set-TextValue "E50" "0.10%"
